$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "didn't ask",
    "asked",
    "didn't ask",
    "no need to ask",
    "no need to ask",
    "no need to ask",
    "didn't ask",
    "asked"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $values[$i]
}

$ws.Range("B14").Select()
